$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 11:50"

# Apply the updated COVID-19 country data (new totals + reordering from the refreshed source)
# Row 4
$ws.Range("B4").Value = 142737
$ws.Range("C4").Value = 277
$ws.Range("E4").Value = 135686

# Row 6
$ws.Range("A6").Value = "España"
$ws.Range("B6").Value = 85195
$ws.Range("C6").Value = 5085
$ws.Range("D6").Value = 16780
$ws.Range("E6").Value = 61075
$ws.Range("F6").Value = 4165
$ws.Range("G6").Value = 537
$ws.Range("H6").Value = 7340

# Row 7
$ws.Range("A7").Value = "China"
$ws.Range("B7").Value = 81470
$ws.Range("C7").Value = 31
$ws.Range("D7").Value = 75700
$ws.Range("E7").Value = 2466
$ws.Range("F7").Value = 633
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 3304

# Row 12
$ws.Range("B12").Value = 15069
$ws.Range("C12").Value = 240
$ws.Range("E12").Value = 12934
$ws.Range("G12").Value = 12
$ws.Range("H12").Value = 312

# Row 26
$ws.Range("D26").Value = 479
$ws.Range("E26").Value = 2110

# Row 34
$ws.Range("A34").Value = "Rusia"
$ws.Range("B34").Value = 1836
$ws.Range("C34").Value = 302
$ws.Range("D34").Value = 66
$ws.Range("E34").Value = 1761
$ws.Range("F34").Value = 8
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 9

# Row 35
$ws.Range("A35").Value = "Rumania"
$ws.Range("B35").Value = 1815
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 206
$ws.Range("E35").Value = 1565
$ws.Range("F35").Value = 31
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 44

# Row 36
$ws.Range("A36").Value = "Pakistan"
$ws.Range("B36").Value = 1625
$ws.Range("C36").Value = 28
$ws.Range("D36").Value = 29
$ws.Range("E36").Value = 1578
$ws.Range("F36").Value = 11
$ws.Range("G36").Value = 4
$ws.Range("H36").Value = 18

# Row 37
$ws.Range("A37").Value = "Filipinas"
$ws.Range("B37").Value = 1546
$ws.Range("C37").Value = 128
$ws.Range("D37").Value = 42
$ws.Range("E37").Value = 1426
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 7
$ws.Range("H37").Value = 78

# Row 52
$ws.Range("A52").Value = "Eslovenia"
$ws.Range("B52").Value = 756
$ws.Range("C52").Value = 26
$ws.Range("D52").Value = 10
$ws.Range("E52").Value = 735
$ws.Range("F52").Value = 28
$ws.Range("H52").Value = 11

# Row 53
$ws.Range("A53").Value = "Serbia"
$ws.Range("B53").Value = 741
$ws.Range("D53").Value = 42
$ws.Range("E53").Value = 686
$ws.Range("F53").Value = 25
$ws.Range("H53").Value = 13

# Row 68
$ws.Range("F68").Value = 15

# Row 71
$ws.Range("B71").Value = 446
$ws.Range("C71").Value = 8
$ws.Range("D71").Value = 32
$ws.Range("E71").Value = 403
$ws.Range("F71").Value = 3
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 11

# Row 82
$ws.Range("A82").Value = "Kuwait"
$ws.Range("B82").Value = 266
$ws.Range("C82").Value = 11
$ws.Range("D82").Value = 72
$ws.Range("E82").Value = 194
$ws.Range("F82").Value = 13
$ws.Range("H82").Value = 0

# Row 83
$ws.Range("A83").Value = "Moldavia"
$ws.Range("B83").Value = 263
$ws.Range("D83").Value = 13
$ws.Range("E83").Value = 248
$ws.Range("F83").Value = 33
$ws.Range("H83").Value = 2

# Row 84
$ws.Range("A84").Value = "Republica de Macedonia"
$ws.Range("D84").Value = 3
$ws.Range("E84").Value = 250
$ws.Range("F84").Value = 1
$ws.Range("H84").Value = 6

# Row 85
$ws.Range("A85").Value = "Jordania"
$ws.Range("B85").Value = 259
$ws.Range("D85").Value = 18
$ws.Range("E85").Value = 238
$ws.Range("F85").Value = 3
$ws.Range("H85").Value = 3

# Row 103
$ws.Range("B103").Value = 127
$ws.Range("C103").Value = 1
$ws.Range("D103").Value = 38
$ws.Range("E103").Value = 88
$ws.Range("F103").Value = 3

# Row 110
$ws.Range("A110").Value = "Camboya"
$ws.Range("B110").Value = 107
$ws.Range("C110").Value = 4
$ws.Range("D110").Value = 21
$ws.Range("E110").Value = 86
$ws.Range("F110").Value = 1
$ws.Range("H110").Value = 0

# Row 111
$ws.Range("A111").Value = "Guadalupe"
$ws.Range("B111").Value = 106
$ws.Range("D111").Value = 17
$ws.Range("E111").Value = 85
$ws.Range("F111").Value = 10
$ws.Range("H111").Value = 4

# Row 133
$ws.Range("B133").Value = 38
$ws.Range("C133").Value = 1
$ws.Range("E133").Value = 28

# Row 143
$ws.Range("A143").Value = "Etiopia"
$ws.Range("B143").Value = 23
$ws.Range("C143").Value = 2
$ws.Range("D143").Value = 1
$ws.Range("E143").Value = 22

# Row 144
$ws.Range("A144").Value = "Bermudas"
$ws.Range("B144").Value = 22
$ws.Range("C144").Value = 2
$ws.Range("D144").Value = 2
